# Lecture partielle de l'EDT M1 MIAGE.
# Shift the schedule dates forward by 3 years (2023 -> 2026, same month/day),
# which also changes the French weekday-name labels in column B, and fixes
# the stray time label in D20 to match the other "15:45" entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: date serial numbers (2023 -> 2026)
$ws.Range("A2").Value  = 46049
$ws.Range("A4").Value  = 46051
$ws.Range("A7").Value  = 46055
$ws.Range("A9").Value  = 46057
$ws.Range("A12").Value = 46093
$ws.Range("A16").Value = 46099
$ws.Range("A19").Value = 46107
$ws.Range("A22").Value = 46114

# Column B: French weekday names that correspond to the new dates above
$ws.Range("B2").Value  = "mardi"
$ws.Range("B4").Value  = "jeudi"
$ws.Range("B7").Value  = "lundi"
$ws.Range("B9").Value  = "mercredi"
$ws.Range("B12").Value = "jeudi"
$ws.Range("B16").Value = "mercredi"
$ws.Range("B19").Value = "jeudi"
$ws.Range("B22").Value = "jeudi"

# D20 held a stray "16:45" time that should read "15:45" like the other rows
$ws.Range("D20").Value = "15:45"
